# Weekly Plan Update 6 April 2018
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("1 April 2018")
$ws3 = $wb.Worksheets.Item("8 April 2018")

# --- Row 3 (E3/F3 need the Percent / Date formatting already used on row 2) ---
$ws2.Range("E2").Copy()
$ws2.Range("E3").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("D3").Copy()
$ws2.Range("F3").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 43196
$ws2.Range("G3").Value = 2

# --- Row 4 ---
$ws2.Range("E2").Copy()
$ws2.Range("E4").PasteSpecial(-4122)
$ws2.Range("D4").Copy()
$ws2.Range("F4").PasteSpecial(-4122)
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 43196
$ws2.Range("G4").Value = 0.25

# --- Row 7: actual effort updated from 2 to 3 ---
$ws2.Range("G7").Value = 3

# --- Row 8 ---
$ws2.Range("E2").Copy()
$ws2.Range("E8").PasteSpecial(-4122)
$ws2.Range("D8").Copy()
$ws2.Range("F8").PasteSpecial(-4122)
$ws2.Range("E8").Value = 1
$ws2.Range("F8").Value = 43196
$ws2.Range("G8").Value = 0.25

$excel.CutCopyMode = 0

# --- Selections left by the editor when they saved ---
$ws2.Range("G13").Select()
$ws3.Range("C5").Select()

# Make sure the sheet that was active/visible when saved matches the source
$ws2.Activate()
